# plantilla_clientes.xlsx edit: drop SECTOR ECONOMICO / AREA DE TRABAJO / LECTOESCRITURA
# columns, rename the birth-date header, and store the birth date as free text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# 1. Rename header K1: FECHA-NACIMIENTO -> FECHA DE NACIMIENTO
$ws.Range("K1").Value = "FECHA DE NACIMIENTO"

# 2. K2 becomes free text instead of a real date value.
$ws.Range("K1:K1048576").NumberFormat = "@"
$ws.Range("K2").Value = "01/05/2004"

# 3. Update DOCUMENTO value for the sample row.
$ws.Range("F2").Value = 2222132

# 4. Drop the columns that are no longer part of the template:
#    P = SECTOR ECONOMICO, Q = AREA DE TRABAJO (adjacent, delete together),
#    then (old) S = LECTOESCRITURA, which after the first delete is column Q.
$ws.Columns("P:Q").Delete()
$ws.Columns("Q:Q").Delete()

# 5. Fix up the sheet view to match the new, narrower layout.
$ws.Range("I1").Select()
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("P1:P1048576").Select()
